$d = $word.ActiveDocument

# Helper: find literal text in the whole document and replace the matched
# range's contents via Range.InsertXML so that sibling runs (in particular
# the leading zero-width <w:r/> that several paragraphs use) are left alone
# - a plain Find.Execute(...,Replace:=wdReplaceAll) or Range.Text assignment
# normalizes/merges the paragraph's runs and drops that empty run, and it
# also lets AutoCorrect mangle straight quotes into curly ones.
function Replace-Xml {
    param(
        [string]$OldText,
        [string]$RunXml
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $OldText"
    }

    # Re-wrap the Find hit's coordinates into a plain Range. Using the
    # Find-produced range object directly makes InsertXML *append* after the
    # matched text instead of replacing it; a freshly constructed Range with
    # the same Start/End behaves like a normal replace and leaves sibling
    # runs (e.g. a leading zero-width <w:r/>) untouched.
    $target = $d.Range($rng.Start, $rng.End)

    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $RunXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($pkg)
}

function Esc {
    param([string]$s)
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# 1. Main title (Heading1)
Replace-Xml "Play Money Mouse Slot Free | Chinese New Year Theme" `
    ('<w:r><w:t xml:space="preserve">' + (Esc "Play Money Mouse Free - Review of Pragmatic Play's Slot Game") + '</w:t></w:r>')

# 2. "What we like" bullet list - process from bottom to top so we don't
#    clobber text we still need to match further down.
Replace-Xml "Celebration of Chinese culture." `
    ('<w:r><w:t xml:space="preserve">' + (Esc "Celebrates Chinese New Year and Chinese culture") + '</w:t></w:r>')

Replace-Xml "Reference to Chinese New year theme." `
    ('<w:r><w:t xml:space="preserve">' + (Esc "Three jackpots of different values") + '</w:t></w:r>')

Replace-Xml "Three jackpots of different values." `
    ('<w:r><w:t xml:space="preserve">' + (Esc "Colorful and funny graphics") + '</w:t></w:r>')

Replace-Xml "Colorful and funny graphics." `
    ('<w:r><w:t xml:space="preserve">' + (Esc "Gameplay with a Wild symbol and bonus mode") + '</w:t></w:r>')

# 3. "What we don't like" bullet list
Replace-Xml "Difficult to hit jackpots." `
    ('<w:r><w:t xml:space="preserve">' + (Esc "Difficulty in obtaining the jackpots") + '</w:t></w:r>')

Replace-Xml "Limited symbol variety." `
    ('<w:r><w:t xml:space="preserve">' + (Esc "Highly used oriental theme in the online slot game industry") + '</w:t></w:r>')

# 4. Bold text (near bottom, duplicate of title)
Replace-Xml "Play Money Mouse Slot Free | Chinese New Year Theme" `
    ('<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">' + (Esc "Play Money Mouse Free - Review of Pragmatic Play's Slot Game") + '</w:t></w:r>')

# 5. Italic text (meta description)
Replace-Xml "Money Mouse by Pragmatic Play celebrates Chinese New Year with colorful graphics, 3 jackpots, and references to Chinese culture. Play for free." `
    ('<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">' + (Esc "Review of Money Mouse by Pragmatic Play - gameplay, graphics, jackpots, and Chinese New Year theme. Play free!") + '</w:t></w:r>')
